$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.758.59"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "'2.101.37"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'227.61"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "'62.23"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").Value = "'0.0840"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "'15.73"
$ws.Range("E12").Value = "  +6.16%  "
$ws.Range("D13").Value = "'2.412.99"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'0.807"
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("D16").Value = "'5.54"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "'2.101.61"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "'38.745.64"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "'71.87"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "'0.0₃0839"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'227.52"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "'9.62"
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("D27").Value = "'172.20"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +4.04%  "
$ws.Range("E29").Value = "  +4.84%  "
$ws.Range("D30").Value = "'19.33"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("E31").Value = "  +10.39%  "
$ws.Range("D32").Value = "'0.120"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'4.54"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("D34").Value = "'4.76"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").Value = "'7.00"
$ws.Range("E35").Value = "  +6.72%  "
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "'3.59"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'18.10"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "'102.75"
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("E42").Value = "  +3.96%  "
$ws.Range("D43").Value = "'1.525.80"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("E44").Value = "  +7.25%  "
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "'0.0911"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "'4.15"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "'1.05"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "'2.298.88"
$ws.Range("E51").Value = "  +0.25%  "
